$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (B2:O2) gets the same values as row 3 (B3:O3)
$values = @(2,2,2,1,1,1,2,1,2,2,2,1,2,2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $values[$i]
}
